$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 1.918831666666666
$ws.Range("N2").Value = 5.756494999999999
$ws.Range("O2").Value = 0.2136389939136679
$ws.Range("P2").Value = 0.2136389939136679
$ws.Range("Q2").Value = 68.57683966869222
$ws.Range("R2").Value = 617.1915570182299
$ws.Range("S2").Value = 0.004165382293258938
$ws.Range("T2").Value = 0.004165382293258936
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.3598037466296114
$ws.Range("P3").Value = 0.3598037466296114
$ws.Range("Q3").Value = 115.494851350894
$ws.Range("R3").Value = 1039.453662158046
$ws.Range("S3").Value = 0.007015199462439172
$ws.Range("T3").Value = 0.00701519946243917
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 3.26469
$ws.Range("N4").Value = 9.794070000000001
$ws.Range("O4").Value = 0.3634842488562984
$ws.Range("P4").Value = 0.3634842488562984
$ws.Range("Q4").Value = 116.67627055942
$ws.Range("R4").Value = 1050.08643503478
$ws.Range("S4").Value = 0.007086959296748902
$ws.Range("T4").Value = 0.007086959296748901
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("M5").Value = 0.5665
$ws.Range("N5").Value = 1.6995
$ws.Range("O5").Value = 0.06307301060042241
$ws.Range("P5").Value = 0.06307301060042241
$ws.Range("Q5").Value = 20.24605928033333
$ws.Range("R5").Value = 182.214533523
$ws.Range("S5").Value = 0.001229753036768652
$ws.Range("T5").Value = 0.001229753036768652
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 1.918831666666666
$ws.Range("N6").Value = 5.756494999999999
$ws.Range("O6").Value = 0.2136389939136679
$ws.Range("P6").Value = 0.2136389939136679
$ws.Range("Q6").Value = 3241.461815793377
$ws.Range("R6").Value = 29173.1563421404
$ws.Range("S6").Value = 0.1968875748286314
$ws.Range("T6").Value = 0.1968875748286314
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.3598037466296114
$ws.Range("P7").Value = 0.3598037466296114
$ws.Range("R7").Value = 49132.46762974007
$ws.Range("S7").Value = 0.3315915591550976
$ws.Range("T7").Value = 0.3315915591550976
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 3.26469
$ws.Range("N8").Value = 9.794070000000001
$ws.Range("O8").Value = 0.3634842488562984
$ws.Range("P8").Value = 0.3634842488562984
$ws.Range("Q8").Value = 5515.0059065816
$ws.Range("R8").Value = 49635.05315923441
$ws.Range("S8").Value = 0.3349834734507464
$ws.Range("T8").Value = 0.3349834734507464
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("M9").Value = 0.5665
$ws.Range("N9").Value = 1.6995
$ws.Range("O9").Value = 0.06307301060042241
$ws.Range("P9").Value = 0.06307301060042241
$ws.Range("Q9").Value = 956.9823922266667
$ws.Range("R9").Value = 8612.841530039999
$ws.Range("S9").Value = 0.05812746009876826
$ws.Range("T9").Value = 0.05812746009876826
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 1.918831666666666
$ws.Range("N10").Value = 5.756494999999999
$ws.Range("O10").Value = 0.2136389939136679
$ws.Range("P10").Value = 0.2136389939136679
$ws.Range("Q10").Value = 179.682377212245
$ws.Range("R10").Value = 1617.141394910205
$ws.Range("S10").Value = 0.01091397323158143
$ws.Range("T10").Value = 0.01091397323158143
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.3598037466296114
$ws.Range("P11").Value = 0.3598037466296114
$ws.Range("Q11").Value = 302.615132846049
$ws.Range("R11").Value = 2723.536195614441
$ws.Range("S11").Value = 0.0183809537173029
$ws.Range("T11").Value = 0.0183809537173029
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 3.26469
$ws.Range("N12").Value = 9.794070000000001
$ws.Range("O12").Value = 0.3634842488562984
$ws.Range("P12").Value = 0.3634842488562984
$ws.Range("Q12").Value = 305.71064166357
$ws.Range("R12").Value = 2751.395774972131
$ws.Range("S12").Value = 0.01856897605369843
$ws.Range("T12").Value = 0.01856897605369843
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("M13").Value = 0.5665
$ws.Range("N13").Value = 1.6995
$ws.Range("O13").Value = 0.06307301060042241
$ws.Range("P13").Value = 0.06307301060042241
$ws.Range("Q13").Value = 53.0479397745
$ws.Range("R13").Value = 477.4314579705
$ws.Range("S13").Value = 0.003222151240828428
$ws.Range("T13").Value = 0.003222151240828428
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 1.918831666666666
$ws.Range("N14").Value = 5.756494999999999
$ws.Range("O14").Value = 0.2136389939136679
$ws.Range("P14").Value = 0.2136389939136679
$ws.Range("Q14").Value = 27.52804583363278
$ws.Range("R14").Value = 247.752412502695
$ws.Range("S14").Value = 0.001672063560196155
$ws.Range("T14").Value = 0.001672063560196155
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.3598037466296114
$ws.Range("P15").Value = 0.3598037466296114
$ws.Range("Q15").Value = 46.361826775571
$ws.Range("R15").Value = 417.256440980139
$ws.Range("S15").Value = 0.002816034294771756
$ws.Range("T15").Value = 0.002816034294771756
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 3.26469
$ws.Range("N16").Value = 9.794070000000001
$ws.Range("O16").Value = 0.3634842488562984
$ws.Range("P16").Value = 0.3634842488562984
$ws.Range("Q16").Value = 46.83607088303
$ws.Range("R16").Value = 421.5246379472701
$ws.Range("S16").Value = 0.002844840055104774
$ws.Range("T16").Value = 0.002844840055104774
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("M17").Value = 0.5665
$ws.Range("N17").Value = 1.6995
$ws.Range("O17").Value = 0.06307301060042241
$ws.Range("P17").Value = 0.06307301060042241
$ws.Range("Q17").Value = 53.0479397745
$ws.Range("R17").Value = 477.4314579705
$ws.Range("S17").Value = 0.003222151240828428
$ws.Range("T17").Value = 0.003222151240828428
